{"js": "// The bold phrase \"DOCX, DOC, PDF, HTML, XPS, R\" / \"TF and TXT\" was split into\n// two runs by a \"_GoBack\" bookmark. Re-write it as a single, unbroken run\n// (search treats the body text as one continuous stream, so it matches across\n// the run/bookmark boundary) and drop the now-unused bookmark.\nconst body = context.document.body;\n\nconst results = body.search(\"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n    results.items[0].insertText(\n        \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\",\n        Word.InsertLocation.replace\n    );\n    await context.sync();\n}\n\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The bold run \"DOCX, DOC, PDF, HTML, XPS, R\" was split from \"TF and TXT\" by a\n# \"_GoBack\" bookmark. Replace the full bold phrase (Find/Replace treats the\n# content as one continuous text stream, so it matches across the run/bookmark\n# boundary) with the same text as a single, unbroken run and let Word drop the\n# now-empty bookmark markers.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Remove the now-orphaned \"_GoBack\" bookmark, if Word kept it around.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
